$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '51.685.48'),
    @('E2', '  -0.92%  '),
    @('D3', '2.924.45'),
    @('E3', '  +1.37%  '),
    @('D4', '0.999'),
    @('E4', '  -0.15%  '),
    @('D5', '354.64'),
    @('E5', '  +0.37%  '),
    @('D6', '109.89'),
    @('E6', '  -1.90%  '),
    @('D7', '0.562'),
    @('E7', '  +0.58%  '),
    @('E8', '  -0.05%  '),
    @('D9', '0.631'),
    @('E9', '  +1.27%  '),
    @('D10', '39.17'),
    @('E10', '  -2.81%  '),
    @('D11', '0.0896'),
    @('E11', '  +4.09%  '),
    @('E12', '  +0.86%  '),
    @('D13', '19.75'),
    @('E13', '  -1.97%  '),
    @('D14', '7.94'),
    @('E14', '  +1.28%  '),
    @('D15', '3.380.06'),
    @('E15', '  +1.26%  '),
    @('D16', '2.929.01'),
    @('E16', '  +1.23%  '),
    @('E17', '  -1.73%  '),
    @('D18', '51.687.40'),
    @('E18', '  -0.94%  '),
    @('D19', '7.59'),
    @('E19', '  -0.65%  '),
    @('D20', '3.28'),
    @('E20', '  -2.22%  '),
    @('D21', '14.19'),
    @('E21', '  +3.17%  '),
    @('D22', '0.0₃0984'),
    @('E22', '  +0.06%  '),
    @('D23', '70.90'),
    @('E23', '  -0.50%  '),
    @('D24', '269.95'),
    @('E24', '  -0.02%  '),
    @('E25', '  +0.78%  '),
    @('E26', '  +11.37%  '),
    @('D27', '27.21'),
    @('E27', '  +2.94%  '),
    @('D28', '7.58'),
    @('E28', '  +19.24%  '),
    @('E29', '  +0.05%  '),
    @('E30', '  +14.28%  '),
    @('D31', '10.62'),
    @('E31', '  +0.21%  '),
    @('D32', '38.23'),
    @('E32', '  -1.30%  '),
    @('D33', '6.15'),
    @('E33', '  +2.83%  '),
    @('D34', '52.27'),
    @('E34', '  -1.69%  '),
    @('B35', 'Toncoin'),
    @('C35', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @('D35', '2.07'),
    @('E35', '  -8.49%  '),
    @('B36', 'VeChain'),
    @('C36', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D36', '0.0441'),
    @('E36', '  -4.64%  '),
    @('B37', 'FirstDigitalUSD'),
    @('C37', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'),
    @('D37', '0.999'),
    @('E37', '  +0.01%  '),
    @('E38', '  -2.48%  '),
    @('D39', '18.33'),
    @('E39', '  -1.79%  '),
    @('D40', '2.02'),
    @('E40', '  -1.00%  '),
    @('D41', '2.74'),
    @('E41', '  +4.34%  '),
    @('E42', '  +0.90%  '),
    @('D43', '23.19'),
    @('E43', '  +2.24%  '),
    @('E44', '  -2.14%  '),
    @('B45', 'Monero'),
    @('C45', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @('D45', '117.01'),
    @('E45', '  -3.92%  '),
    @('B46', 'ApeXProtocol'),
    @('C46', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'),
    @('D46', '2.52'),
    @('E46', '  +2.13%  '),
    @('D47', '3.45'),
    @('E47', '  -3.73%  '),
    @('D48', '2.140.06'),
    @('E48', '  -2.34%  '),
    @('E49', '  -5.75%  '),
    @('E50', '  +2.14%  '),
    @('B51', 'MultiversX'),
    @('C51', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'),
    @('D51', '61.95'),
    @('E51', '  +3.28%  ')
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}
